# Update the "F" column (活动热度/数量 type numeric field) on three sheets
# to reflect refreshed output values from the data source regeneration.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1788
$ws1.Range("F4").Value  = 454
$ws1.Range("F6").Value  = 69
$ws1.Range("F9").Value  = 1734
$ws1.Range("F10").Value = 366
$ws1.Range("F13").Value = 337
$ws1.Range("F14").Value = 681
$ws1.Range("F15").Value = 12807
$ws1.Range("F16").Value = 12802
$ws1.Range("F17").Value = 956
$ws1.Range("F21").Value = 51
$ws1.Range("F22").Value = 565
$ws1.Range("F27").Value = 33
$ws1.Range("F28").Value = 250

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 84
$ws2.Range("F7").Value = 8

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1788
$ws4.Range("F6").Value  = 454
$ws4.Range("F9").Value  = 69
$ws4.Range("F14").Value = 1734
$ws4.Range("F15").Value = 366
$ws4.Range("F18").Value = 337
$ws4.Range("F19").Value = 84
$ws4.Range("F20").Value = 681
$ws4.Range("F21").Value = 12807
$ws4.Range("F22").Value = 12802
$ws4.Range("F23").Value = 956
$ws4.Range("F27").Value = 51
$ws4.Range("F28").Value = 565
$ws4.Range("F30").Value = 8
$ws4.Range("F37").Value = 33
$ws4.Range("F38").Value = 250

$wb.Save()
